$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: bus 632, A-N voltage changes from 0.0000 to 1.0000
$ws.Range("B2").Value = "1.0000 at 0.00"

# Update row 3: bus 645 -> 633, with new voltage readings
$ws.Range("A3").Value = 633
$ws.Range("B3").Value = "0.9969 at -0.07"
$ws.Range("C3").Value = "0.9980 at -120.05"
$ws.Range("D3").Value = "0.9973 at 119.99"

# Update row 4: bus 646 -> 634, with new voltage readings
$ws.Range("A4").Value = 634
$ws.Range("B4").Value = "0.9724 at -0.77"
$ws.Range("C4").Value = "0.9788 at -120.55"
$ws.Range("D4").Value = "0.9782 at 119.50"

# Add new row 5: bus 671
$ws.Range("A5").Value = 671
$ws.Range("B5").Value = "1.0011 at 0.01"
$ws.Range("C5").Value = "0.9980 at -120.01"
$ws.Range("D5").Value = "0.9973 at 119.81"
